$wb = $excel.ActiveWorkbook

# --- Week 13 logging: remove M.Sargent from the RB roster ---
$wsRB = $wb.Worksheets.Item("RB")
$wsRB.Rows.Item(4).Delete()

# Make the RB sheet the active tab, with D12 selected (matches the
# author's final view state after finishing the edit).
$wsRB.Activate() | Out-Null
$wsRB.Range("D12").Select() | Out-Null
